$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '56.253.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +9.31%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.241.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.75%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '401.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.559'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0899'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.758.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.08'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.42%  '
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.68%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.257.39'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '56.169.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.32%  '
$ws.Range("E22").Value = '  +5.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '289.40'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.80%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.26'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.24%  '
$ws.Range("E29").Value = '  +2.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("E31").Value = '  +4.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.32%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '37.18'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.47%  '
$ws.Range("B34").Value = 'VeChain'
$ws.Range("C34").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0498'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.01%  '
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.32'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.01%  '
$ws.Range("E37").Value = '  +7.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.07'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +21.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '139.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.94'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.04'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +10.07%  '
$ws.Range("E43").Value = '  -3.60%  '
$ws.Range("E44").Value = '  +3.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.02'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.73'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("E47").Value = '  +41.22%  '
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.146.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.49%  '
$ws.Range("E50").Value = '  +0.68%  '
$ws.Range("B51").Value = 'BEAM'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0349'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.27%  '
